$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing column B (old B:E shifts to E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# Set new header row text for the inserted columns
$ws.Range("B1").Value = "Absorbace Filtrate"
$ws.Range("C1").Value = "Absorbace S1"
$ws.Range("D1").Value = "Absorbace P1"

# Fill in the new data columns (B = Absorbace Filtrate, C = Absorbace S1, D = Absorbace P1)
$newData = @(
    @("NaN", "NaN", "1.512"),
    @("NaN", "NaN", "NaN"),
    @("NaN", "2.57", "NaN"),
    @("2.421", "2.785", "NaN"),
    @("NaN", "3.099", "NaN"),
    @("3.925", "3.742", "NaN"),
    @("4.243", "4.536", "NaN"),
    @("NaN", "4.89", "4.469"),
    @("NaN", "3.999", "NaN"),
    @("4.428", "4.17", "4.416"),
    @("4.128", "4.621", "4.603"),
    @("4.463", "3.901", "3.969"),
    @("4.157", "4.057", "4.359"),
    @("3.878", "NaN", "4.195"),
    @("4.26", "4.062", "4.259"),
    @("3.895", "3.989", "3.972"),
    @("3.904", "3.957", "3.882"),
    @("3.909", "4.002", "3.947"),
    @("4.008", "3.821", "3.961"),
    @("3.782", "3.824", "3.826"),
    @("3.707", "3.824", "3.765"),
    @("3.694", "3.679", "3.718"),
    @("3.697", "3.708", "3.619"),
    @("3.564", "3.702", "3.644"),
    @("3.571", "3.759", "3.595"),
    @("3.559", "3.502", "3.601"),
    @("3.637", "3.693", "3.591"),
    @("3.538", "3.575", "3.511"),
    @("3.54", "3.464", "3.479"),
    @("3.451", "3.454", "3.43"),
    @("3.5", "3.346", "3.432"),
    @("3.391", "3.428", "3.434"),
    @("3.529", "3.472", "3.442"),
    @("3.417", "3.402", "3.529"),
    @("3.437", "3.302", "3.533"),
    @("3.479", "3.55", "3.417"),
    @("3.461", "3.434", "3.377"),
    @("3.548", "3.453", "3.473"),
    @("3.37", "3.486", "3.385"),
    @("3.499", "3.446", "3.266"),
    @("3.433", "3.382", "2.877")
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $i + 2
    $rowVals = $newData[$i]
    $ws.Cells.Item($r, 2).Value = $rowVals[0]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
}

